# Расчет анкерных болтов - обновление входных данных расчета.
# Единственный лист книги - "Фланец".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Фланец")

# Выбор другого болта в выпадающем списке H3 (список K3:K10): было М48, стало М42.
$ws.Range("H3").Value = "М42"

# Новое значение продольной силы (A4): было 1700, стало 1670.
$ws.Range("A4").Value = 1670

# Новые исходные данные для расчета напряжений (E18, E19).
$ws.Range("E18").Value = 3612.88
$ws.Range("E19").Value = 185.82
